$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell B3 value to "-"
$ws.Range("B3").Value = "-"

# Update the active selection to B3
$ws.Range("B3").Select()
